$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$ws.Range("D2").Value = "63.577.14"
$ws.Range("D3").Value = "3.089.50"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "587.72"
$ws.Range("E5").Value = "  -0.92%  "
Set-TextValue $ws.Range("D6") "152.69"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.552"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").Value = "3.083.78"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  -2.44%  "
Set-TextValue $ws.Range("D11") "5.87"
$ws.Range("E11").Value = "  -1.15%  "
Set-TextValue $ws.Range("D12") "0.461"
$ws.Range("E12").Value = "  -0.76%  "
Set-TextValue $ws.Range("D13") "37.58"
$ws.Range("E13").Value = "  +0.31%  "
Set-TextValue $ws.Range("D14") "0.0000242"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "3.600.91"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("E16").Value = "  -1.95%  "
Set-TextValue $ws.Range("D17") "7.17"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "63.578.37"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "3.086.44"
$ws.Range("E19").Value = "  -2.19%  "
Set-TextValue $ws.Range("D20") "476.43"
$ws.Range("E20").Value = "  +1.41%  "
Set-TextValue $ws.Range("D21") "14.67"
$ws.Range("E21").Value = "  +1.27%  "
Set-TextValue $ws.Range("D22") "0.718"
$ws.Range("E22").Value = "  -2.55%  "
Set-TextValue $ws.Range("D23") "7.56"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D24") "13.14"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D25") "2.35"
$ws.Range("E25").Value = "  +0.49%  "
Set-TextValue $ws.Range("D26") "81.91"
$ws.Range("E26").Value = "  +0.26%  "
Set-TextValue $ws.Range("D27") "0.998"
$ws.Range("E27").Value = "  -0.07%  "
Set-TextValue $ws.Range("D28") "9.69"
$ws.Range("E28").Value = "  +0.32%  "
Set-TextValue $ws.Range("D29") "2.68"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("E30").Value = "  +0.02%  "
Set-TextValue $ws.Range("D31") "7.25"
$ws.Range("E31").Value = "  -3.30%  "
Set-TextValue $ws.Range("D32") "2.20"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("E33").Value = "  +2.95%  "
Set-TextValue $ws.Range("D34") "27.49"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  -1.64%  "
Set-TextValue $ws.Range("D37") "3.40"
$ws.Range("E37").Value = "  +4.14%  "
Set-TextValue $ws.Range("D38") "6.12"
$ws.Range("E38").Value = "  -2.27%  "
Set-TextValue $ws.Range("D39") "2.23"
$ws.Range("E39").Value = "  -5.10%  "
Set-TextValue $ws.Range("D40") "9.33"
$ws.Range("E40").Value = "  +0.58%  "
Set-TextValue $ws.Range("D41") "50.62"
Set-TextValue $ws.Range("D42") "445.24"
$ws.Range("E42").Value = "  -2.92%  "
Set-TextValue $ws.Range("D43") "0.284"
$ws.Range("E43").Value = "  -4.48%  "
Set-TextValue $ws.Range("D44") "0.0365"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "2.825.47"
$ws.Range("E45").Value = "  -4.03%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D46") "0.109"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D47") "38.88"
$ws.Range("E47").Value = "  -2.68%  "
Set-TextValue $ws.Range("D48") "129.96"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D49") "0.999"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D50") "25.25"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("E51").Value = "  +0.59%  "
